$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" everywhere it appears ---

# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: Status column (C), rows 2-4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: Status column (C), rows 2-4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Column width changes: 17.2159881591797 -> 13.4101845877511 (chars) ---
# The COM ColumnWidth property uses a different unit than the stored XML
# character width (XML width = ColumnWidth + 5/6, rounded to the engine's
# pixel grid), so 12.42 is the COM value that lands closest to the target.

$wsOverview.Columns.Item(5).ColumnWidth = 12.42
$wsOverview.Columns.Item(6).ColumnWidth = 12.42

$wsZhCn.Columns.Item(3).ColumnWidth = 12.42

$wsDeDe.Columns.Item(3).ColumnWidth = 12.42
